# Update "想去人数" (want-to-go count) figures in the F column for both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same rows).
# Values taken from the latest site scrape (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 20
$ws1.Range("F3").Value  = 1795
$ws1.Range("F5").Value  = 775
$ws1.Range("F7").Value  = 108
$ws1.Range("F8").Value  = 35
$ws1.Range("F15").Value = 4187
$ws1.Range("F20").Value = 969
$ws1.Range("F21").Value = 1337
$ws1.Range("F23").Value = 38
$ws1.Range("F24").Value = 42
$ws1.Range("F25").Value = 1932
$ws1.Range("F26").Value = 63
$ws1.Range("F28").Value = 88

# --- Sheet "全部类型" (same rows, F21 differs slightly from "展览") ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 20
$ws4.Range("F3").Value  = 1795
$ws4.Range("F5").Value  = 775
$ws4.Range("F7").Value  = 108
$ws4.Range("F8").Value  = 35
$ws4.Range("F15").Value = 4187
$ws4.Range("F20").Value = 969
$ws4.Range("F21").Value = 1338
$ws4.Range("F23").Value = 38
$ws4.Range("F24").Value = 42
$ws4.Range("F25").Value = 1932
$ws4.Range("F26").Value = 63
$ws4.Range("F28").Value = 88
